$d = $word.ActiveDocument

function Get-ParagraphIndexContaining {
    param([string]$needle)
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Add the new "6.El caso de uso finaliza. " paragraph right after
#    "5.El sistema registra la nueva actividad." in the "Flujo
#    basico" cell, carrying the _GoBack bookmark (this is where the
#    author's cursor ended up after the edit).
# ------------------------------------------------------------------
$idx5 = Get-ParagraphIndexContaining "5.El sistema registra la nueva actividad."
if ($idx5 -lt 0) {
    throw "Could not locate paragraph '5.El sistema registra la nueva actividad.'"
}

$para5Range = $d.Paragraphs.Item($idx5).Range
# Position right before this paragraph's end-of-paragraph mark so the
# inserted paragraph lands after it, inside the same table cell,
# without touching paragraph 5's own text.
$insertPoint = $d.Range($para5Range.End - 1, $para5Range.End - 1)

$newParaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">6.El caso de uso finaliza. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($newParaXml)

# ------------------------------------------------------------------
# 2) Remove the old _GoBack bookmark that used to sit in the
#    "2.2 El sistema no soporta el peso del archivo cargado"
#    paragraph - a document only keeps one _GoBack bookmark, and it
#    now belongs to the newly-added paragraph above.
# ------------------------------------------------------------------
$idx22 = Get-ParagraphIndexContaining "2.2 El sistema no soporta el peso del archivo cargado"
if ($idx22 -lt 0) {
    throw "Could not locate paragraph '2.2 El sistema no soporta el peso del archivo cargado'"
}

$para22Range = $d.Paragraphs.Item($idx22).Range

$replacementXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>2.2 El sistema no soporta el peso del archivo cargado</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$para22Range.InsertXML($replacementXml)

Write-Output "done"
